# Auto-generated edit script: updates market-price-derived columns (H-N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets per the scheduled runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 1791.775
$ws.Cells.Item(17, 10).Value = 1844.3823
$ws.Cells.Item(17, 12).Value = 5533.1469
$ws.Cells.Item(17, 14).Value = -5869.1469
$ws.Cells.Item(28, 8).Value = 417
$ws.Cells.Item(28, 9).Value = 463.33334
$ws.Cells.Item(28, 10).Value = 347.5
$ws.Cells.Item(28, 11).Value = 463.33334
$ws.Cells.Item(28, 12).Value = 347.5
$ws.Cells.Item(28, 13).Value = 21.66665999999998
$ws.Cells.Item(28, 14).Value = -1317.5
$ws.Cells.Item(62, 8).Value = 6616.125
$ws.Cells.Item(62, 9).Value = 5825.1665
$ws.Cells.Item(62, 11).Value = 5825.1665
$ws.Cells.Item(62, 13).Value = -5201.1665
$ws.Cells.Item(65, 8).Value = 6616.125
$ws.Cells.Item(65, 9).Value = 5825.1665
$ws.Cells.Item(65, 11).Value = 29125.8325
$ws.Cells.Item(65, 13).Value = -26005.8325
$ws.Cells.Item(92, 8).Value = 409.8125
$ws.Cells.Item(92, 9).Value = 414.48
$ws.Cells.Item(92, 11).Value = 414.48
$ws.Cells.Item(92, 13).Value = 833.52
$ws.Cells.Item(106, 8).Value = 10472.941
$ws.Cells.Item(106, 9).Value = 4717.857
$ws.Cells.Item(106, 11).Value = 4717.857
$ws.Cells.Item(106, 13).Value = -4086.857
$ws.Cells.Item(121, 8).Value = 52262
$ws.Cells.Item(121, 10).Value = 52262
$ws.Cells.Item(121, 12).Value = 156786
$ws.Cells.Item(121, 14).Value = -160280
$ws.Cells.Item(130, 8).Value = 93814.336
$ws.Cells.Item(130, 10).Value = 93814.336
$ws.Cells.Item(130, 12).Value = 93814.336
$ws.Cells.Item(130, 14).Value = -103854.336
$ws.Cells.Item(134, 8).Value = 112572.43
$ws.Cells.Item(134, 10).Value = 107424.305
$ws.Cells.Item(134, 12).Value = 107424.305
$ws.Cells.Item(134, 14).Value = -117564.305
$ws.Cells.Item(137, 8).Value = 3302
$ws.Cells.Item(137, 9).Value = 2604.45
$ws.Cells.Item(137, 10).Value = 6789.75
$ws.Cells.Item(137, 11).Value = 7813.349999999999
$ws.Cells.Item(137, 12).Value = 20369.25
$ws.Cells.Item(137, 13).Value = -5263.349999999999
$ws.Cells.Item(137, 14).Value = -25469.25
$ws.Cells.Item(138, 8).Value = 6119.1777
$ws.Cells.Item(138, 10).Value = 5916.15
$ws.Cells.Item(138, 12).Value = 17748.45
$ws.Cells.Item(138, 14).Value = -28028.45

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 12405.487
$ws.Cells.Item(32, 9).Value = 12405.487
$ws.Cells.Item(32, 11).Value = 12405.487
$ws.Cells.Item(32, 13).Value = -12118.487
$ws.Cells.Item(61, 8).Value = 1586.8572
$ws.Cells.Item(61, 9).Value = 1584.6666
$ws.Cells.Item(61, 10).Value = 1600
$ws.Cells.Item(61, 11).Value = 1584.6666
$ws.Cells.Item(61, 12).Value = 1600
$ws.Cells.Item(61, 13).Value = -1372.6666
$ws.Cells.Item(61, 14).Value = -2024
$ws.Cells.Item(74, 8).Value = 2587.889
$ws.Cells.Item(74, 9).Value = 2298.5
$ws.Cells.Item(74, 10).Value = 3166.6667
$ws.Cells.Item(74, 11).Value = 2298.5
$ws.Cells.Item(74, 12).Value = 3166.6667
$ws.Cells.Item(74, 13).Value = -1424.5
$ws.Cells.Item(74, 14).Value = -4914.6667
$ws.Cells.Item(77, 8).Value = 2587.889
$ws.Cells.Item(77, 9).Value = 2298.5
$ws.Cells.Item(77, 10).Value = 3166.6667
$ws.Cells.Item(77, 11).Value = 11492.5
$ws.Cells.Item(77, 12).Value = 15833.3335
$ws.Cells.Item(77, 13).Value = -7124.5
$ws.Cells.Item(77, 14).Value = -24569.3335
$ws.Cells.Item(102, 8).Value = 1326.2858
$ws.Cells.Item(102, 9).Value = 1369.909
$ws.Cells.Item(102, 11).Value = 1369.909
$ws.Cells.Item(102, 13).Value = 252.0909999999999
$ws.Cells.Item(110, 8).Value = 4667
$ws.Cells.Item(110, 9).Value = 3699.2
$ws.Cells.Item(110, 11).Value = 3699.2
$ws.Cells.Item(110, 13).Value = -1654.2
$ws.Cells.Item(132, 8).Value = 3323
$ws.Cells.Item(132, 9).Value = 3376.0967
$ws.Cells.Item(132, 10).Value = 2500
$ws.Cells.Item(132, 11).Value = 10128.2901
$ws.Cells.Item(132, 12).Value = 7500
$ws.Cells.Item(132, 13).Value = -7598.2901
$ws.Cells.Item(132, 14).Value = -12560
$ws.Cells.Item(136, 8).Value = 1586.8572
$ws.Cells.Item(136, 9).Value = 1584.6666
$ws.Cells.Item(136, 10).Value = 1600
$ws.Cells.Item(136, 11).Value = 4753.9998
$ws.Cells.Item(136, 12).Value = 4800
$ws.Cells.Item(136, 13).Value = -2203.9998
$ws.Cells.Item(136, 14).Value = -9900

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 41678610
$ws.Cells.Item(20, 9).Value = 45466708
$ws.Cells.Item(20, 10).Value = 9498
$ws.Cells.Item(20, 11).Value = 45466708
$ws.Cells.Item(20, 12).Value = 9498
$ws.Cells.Item(20, 13).Value = -45466461
$ws.Cells.Item(20, 14).Value = -9992

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 42415.348
$ws.Cells.Item(31, 9).Value = 61369.707
$ws.Cells.Item(31, 10).Value = 6612.6665
$ws.Cells.Item(31, 11).Value = 61369.707
$ws.Cells.Item(31, 12).Value = 6612.6665
$ws.Cells.Item(31, 13).Value = -61074.707
$ws.Cells.Item(31, 14).Value = -7202.6665
$ws.Cells.Item(34, 8).Value = 42415.348
$ws.Cells.Item(34, 9).Value = 61369.707
$ws.Cells.Item(34, 10).Value = 6612.6665
$ws.Cells.Item(34, 11).Value = 61369.707
$ws.Cells.Item(34, 12).Value = 6612.6665
$ws.Cells.Item(34, 13).Value = -61167.707
$ws.Cells.Item(34, 14).Value = -7016.6665
$ws.Cells.Item(68, 8).Value = 27777.777
$ws.Cells.Item(68, 10).Value = 27500
$ws.Cells.Item(68, 12).Value = 27500
$ws.Cells.Item(68, 14).Value = -28998
$ws.Cells.Item(71, 8).Value = 27777.777
$ws.Cells.Item(71, 10).Value = 27500
$ws.Cells.Item(71, 12).Value = 82500
$ws.Cells.Item(71, 14).Value = -89988
$ws.Cells.Item(99, 8).Value = 26963.166
$ws.Cells.Item(99, 9).Value = 18386.75
$ws.Cells.Item(99, 11).Value = 18386.75
$ws.Cells.Item(99, 13).Value = -16888.75
$ws.Cells.Item(126, 8).Value = 26963.166
$ws.Cells.Item(126, 9).Value = 18386.75
$ws.Cells.Item(126, 11).Value = 55160.25
$ws.Cells.Item(126, 13).Value = -52690.25
$ws.Cells.Item(132, 8).Value = 8939.5
$ws.Cells.Item(132, 9).Value = 9549.375
$ws.Cells.Item(132, 11).Value = 28648.125
$ws.Cells.Item(132, 13).Value = -26118.125
$ws.Cells.Item(134, 8).Value = 5783.9287
$ws.Cells.Item(134, 9).Value = 6697.5
$ws.Cells.Item(134, 11).Value = 20092.5
$ws.Cells.Item(134, 13).Value = -17557.5
$ws.Cells.Item(141, 8).Value = 637297.8
$ws.Cells.Item(141, 10).Value = 702997.75
$ws.Cells.Item(141, 12).Value = 702997.75
$ws.Cells.Item(141, 14).Value = -713357.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1422.4546
$ws.Cells.Item(5, 10).Value = 2566
$ws.Cells.Item(5, 12).Value = 7698
$ws.Cells.Item(5, 14).Value = -7922
$ws.Cells.Item(88, 8).Value = 10359.8
$ws.Cells.Item(88, 9).Value = 2000
$ws.Cells.Item(88, 11).Value = 6000
$ws.Cells.Item(88, 13).Value = -5572
$ws.Cells.Item(91, 8).Value = 10359.8
$ws.Cells.Item(91, 9).Value = 2000
$ws.Cells.Item(91, 11).Value = 6000
$ws.Cells.Item(91, 13).Value = -4518
$ws.Cells.Item(99, 8).Value = 18745.834
$ws.Cells.Item(99, 10).Value = 18745.834
$ws.Cells.Item(99, 12).Value = 56237.50199999999
$ws.Cells.Item(99, 14).Value = -60729.50199999999
$ws.Cells.Item(101, 8).Value = 24099.8
$ws.Cells.Item(101, 10).Value = 24099.8
$ws.Cells.Item(101, 12).Value = 72299.39999999999
$ws.Cells.Item(101, 14).Value = -77167.39999999999
$ws.Cells.Item(105, 8).Value = 13494.5
$ws.Cells.Item(105, 10).Value = 14989
$ws.Cells.Item(105, 12).Value = 44967
$ws.Cells.Item(105, 14).Value = -50209
$ws.Cells.Item(106, 8).Value = 5852.3335
$ws.Cells.Item(106, 10).Value = 6047.421
$ws.Cells.Item(106, 12).Value = 18142.263
$ws.Cells.Item(106, 14).Value = -20034.263
$ws.Cells.Item(112, 8).Value = 253262
$ws.Cells.Item(112, 9).Value = 253262
$ws.Cells.Item(112, 11).Value = 759786
$ws.Cells.Item(112, 13).Value = -758678
$ws.Cells.Item(113, 8).Value = 7166.3335
$ws.Cells.Item(113, 10).Value = 7166.3335
$ws.Cells.Item(113, 12).Value = 21499.0005
$ws.Cells.Item(113, 14).Value = -25839.0005
$ws.Cells.Item(122, 8).Value = 2000
$ws.Cells.Item(122, 9).Value = 1500
$ws.Cells.Item(122, 10).Value = 2500
$ws.Cells.Item(122, 11).Value = 13500
$ws.Cells.Item(122, 12).Value = 22500
$ws.Cells.Item(122, 13).Value = -11050
$ws.Cells.Item(122, 14).Value = -27400
$ws.Cells.Item(135, 8).Value = 1422.4546
$ws.Cells.Item(135, 10).Value = 2566
$ws.Cells.Item(135, 12).Value = 23094
$ws.Cells.Item(135, 14).Value = -28164
$ws.Cells.Item(137, 8).Value = 8339295
$ws.Cells.Item(137, 10).Value = 6910.875
$ws.Cells.Item(137, 12).Value = 20732.625
$ws.Cells.Item(137, 14).Value = -30932.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 2808.7144
$ws.Cells.Item(122, 9).Value = 2826.2
$ws.Cells.Item(122, 11).Value = 8478.599999999999
$ws.Cells.Item(122, 13).Value = -6028.599999999999
$ws.Cells.Item(132, 8).Value = 10001.641
$ws.Cells.Item(132, 9).Value = 9410.826999999999
$ws.Cells.Item(132, 10).Value = 11715
$ws.Cells.Item(132, 11).Value = 28232.481
$ws.Cells.Item(132, 12).Value = 35145
$ws.Cells.Item(132, 13).Value = -25702.481
$ws.Cells.Item(132, 14).Value = -40205

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 47480
$ws.Cells.Item(93, 9).Value = 8893
$ws.Cells.Item(93, 10).Value = 60342.332
$ws.Cells.Item(93, 11).Value = 8893
$ws.Cells.Item(93, 12).Value = 60342.332
$ws.Cells.Item(93, 13).Value = -7645
$ws.Cells.Item(93, 14).Value = -62838.332
$ws.Cells.Item(100, 8).Value = 17181.062
$ws.Cells.Item(100, 9).Value = 6209.75
$ws.Cells.Item(100, 10).Value = 23763.85
$ws.Cells.Item(100, 11).Value = 6209.75
$ws.Cells.Item(100, 12).Value = 23763.85
$ws.Cells.Item(100, 13).Value = -5668.75
$ws.Cells.Item(100, 14).Value = -24845.85
$ws.Cells.Item(132, 8).Value = 5268.385
$ws.Cells.Item(132, 9).Value = 4511.75
$ws.Cells.Item(132, 10).Value = 6479
$ws.Cells.Item(132, 11).Value = 13535.25
$ws.Cells.Item(132, 12).Value = 19437
$ws.Cells.Item(132, 13).Value = -11005.25
$ws.Cells.Item(132, 14).Value = -24497

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 29413472
$ws.Cells.Item(107, 9).Value = 1445.3334
$ws.Cells.Item(107, 10).Value = 100002340
$ws.Cells.Item(107, 11).Value = 4336.0002
$ws.Cells.Item(107, 12).Value = 300007020
$ws.Cells.Item(107, 13).Value = -2416.0002
$ws.Cells.Item(107, 14).Value = -300010860
$ws.Cells.Item(122, 8).Value = 1520
$ws.Cells.Item(122, 9).Value = 1168.5714
$ws.Cells.Item(122, 11).Value = 3505.7142
$ws.Cells.Item(122, 13).Value = -1055.7142
